$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario_info")

# The matsim-run-output dependency is removed; the network CRS now has to be
# supplied directly as a "global" parameter. Insert a new row for it right
# after the existing "sampleSize" row (new row 6), pushing everything else
# down by one.
$ws.Rows.Item(6).Insert()

$ws.Cells.Item(6, 1).Value = "global"
$ws.Cells.Item(6, 3).Value = "EPSG:25832"
$ws.Cells.Item(6, 2).Value = "networkCrs"
$ws.Cells.Item(6, 5).Value = "The coordinate reference system of the network"

# Re-apply the sheet's AutoFilter so its range grows to cover the new row.
$ws.AutoFilterMode = $false
$ws.Range("A1:E25").AutoFilter()

# Keep the workbook-level "_FilterDatabase" hidden name (driven by the
# AutoFilter above) in sync with the new range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=scenario_info!`$A`$1:`$E`$25"
  }
}

$ws.Range("C9").Select()
